# Applies the GPT-4.0 refresh edit to the IESO cover letter.
$d = $word.ActiveDocument

$find0 = @'
From: Austing Dong^lTo: IESO - Corporate Resources
'@
$repl0 = @'
May 22nd, 2023^l^lFrom: Austing Dong^l^lTo: IESO (Independent Electricity System Operator) - Corporate Resources
'@

$find1 = @'
I am writing to express my strong interest in applying for the position of Information Security Student at IESO. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirement.
'@
$repl1 = @'
I am writing to express my strong interest in applying for the position of Information Security Student at IESO. As a University of Waterloo Computer Science undergraduate student, I strongly believe that my technical competencies and academic background are closely in line with the job requirements. I would like to highlight the following for your consideration:
'@

$find2 = @'
I loved computer science as well as developing applications since Middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest, and I did an excellent job in relevant courses in the beginning of my University studies. I found solving business challenges through programming is fascinating because this is the way I feel the sense of accomplishment. Such deep interest in programming and technology has motivated me to deep dive in related fields such as software development, quality assurance and machine learning.
'@
$repl2 = @'
My passion for computer science and application development began in middle school when I earned the gold rank for the USA Computing Olympiad algorithm contest. This enthusiasm continued into my university studies, where I excelled in relevant courses. I find solving business challenges through programming fascinating, as it provides me with a sense of accomplishment. This deep interest in programming and technology has motivated me to explore related fields such as software development, quality assurance, and machine learning.
'@

$find3 = @'
The computer science co-op program at the University of Waterloo offered me a unique opportunity to take on both programming and logical courses. Through working on a massive number of technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including but not limited to object-oriented programming, web application development, artificial intelligence, algorithm design and data abstraction. Such projects can be viewed on my GitHub: https://github.com/AustingDong. One of the biggest projects I led and built was implementing an application that uses AI to extract keywords from articles containing scientific or technical information which helps users quickly locate their desired items based on keywords' weight. This application can be used to quickly get all the important items and keywords from NASA Technical Report Server which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project detail can be found here. Through understanding the project requirements, researching on coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills through interpreting and explaining technical concepts to my teammates while working in a team environment.^l^lI am extremely interested in advancing my career and contributing my skills to IESO. I am particularly interested in the Information Security Student position as it aligns with my passion for technology and my desire to contribute to a secure and reliable energy system. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment.^l^lI am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any question or require additional information.
'@
$repl3 = @'
The computer science co-op program at the University of Waterloo has offered me a unique opportunity to take on both programming and logical courses. Through working on numerous technical projects and assignments from hackathons and school activities, I gained hands-on experience in fields including object-oriented programming, web application development, artificial intelligence, algorithm design, and data abstraction. My projects can be viewed on my GitHub: https://github.com/AustingDong.^l^lOne of the most significant projects I led and built was an application that uses AI to extract keywords from articles containing scientific or technical information. This application helps users quickly locate their desired items based on keyword weight and can be used to efficiently retrieve important items and keywords from the NASA Technical Report Server, which includes hundreds of thousands of items containing scientific and technical information (STI) created or funded by NASA. Project details can be found here.^l^lThrough understanding project requirements, researching coding algorithms, implementing the application with clean code, and incorporating white and black testing in a limited time, I cultivated working experience in developing a technical application in accordance with business requirements. Moreover, I have excellent documentation and communication skills, which I have demonstrated by interpreting and explaining technical concepts to my teammates while working in a team environment. I am confident that my academic and project background has prepared me well, and now is the time I can make valuable contributions to a workplace environment.^l^lI am extremely interested in advancing my career and contributing my skills to IESO. I am willing to answer any preliminary questions you may have. Please feel free to contact me at austingdong@gmail.com or 1-226-789-9109 if you have any questions or require additional information.
'@

$results = @()
$results += $d.Content.Find.Execute($find0, $true, $false, $false, $false, $false, $true, 1, $false, $repl0, 2)
$results += $d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)
$results += $d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)
$results += $d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)

Write-Output ("Find/Replace results: " + ($results -join ", "))
